# MarketBeat rank tracker update - 10th (Jun_27) + 9th (Jun_26) columns added,
# plus two brand-new firms (Benchmark, Evercore ISI) appended as new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room for two new date columns (Jun_26 and Jun_27) by inserting three
#    blank columns before column B. The report generator always writes the
#    newest date column twice side-by-side (a quirk preserved from the source
#    data), so we need B,C,D free for: Jun_27 | Jun_26 | Jun_26 -- the
#    previously-newest columns (old B:E = Jun_17,Jun_15,Jun_13,Jun_10) shift
#    right to E:H.
# ---------------------------------------------------------------------------
$ws.Range("B1:D1").EntireColumn.Insert()

# Give the three new columns the same fixed width as the rest of the date
# columns (8 characters wide, same as the pre-existing C/D/E columns).
$ws.Columns("B:D").ColumnWidth = 7.14

# ---------------------------------------------------------------------------
# 2) Header row: newest date first (column B), then the next-newest, then the
#    rest of the previously existing header values (already shifted right by
#    the column insert above, so only B1:D1 need new values).
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# ---------------------------------------------------------------------------
# 3) Fill the new B:D columns for every existing analyst row with the "UN"
#    (unchanged) placeholder that the tracker uses when nothing happened on
#    that date, then overwrite the handful of cells that have real rating
#    events for Jun_25/Jun_26.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Row 18 (Royal Bank of Canada): 6/25/2018 raised target - recorded under both
# the Jun_26 columns (C & D), matching the green "notable event" highlight
# used elsewhere in this sheet.
$ws.Range("C18").Value = "6/25/2018,Raises Target,Outperform,`$300.00 -> `$310.00"
$ws.Range("D18").Value = "6/25/2018,Raises Target,Outperform,`$300.00 -> `$310.00"
$ws.Range("C18").Interior.Color = 13434828
$ws.Range("D18").Interior.Pattern = -4142

# Row 22 (BidaskClub): 6/26/2018 downgrade - recorded in B, C and D with the
# new highlight color introduced for this update.
$ws.Range("B22").Value = "6/26/2018,Downgrades,Strong-Buy -> Buy,"
$ws.Range("C22").Value = "6/26/2018,Downgrades,Strong-Buy -> Buy,"
$ws.Range("D22").Value = "6/26/2018,Downgrades,Strong-Buy -> Buy,"
$ws.Range("B22").Interior.Color = 13408767
$ws.Range("C22").Interior.Color = 13408767
$ws.Range("D22").Interior.Pattern = -4142

# ---------------------------------------------------------------------------
# 4) Preserve the pre-existing "notable event" green highlight on the two
#    cells that used to live in column E (now shifted to column H).
# ---------------------------------------------------------------------------
$ws.Range("H13").Interior.Color = 13434828
$ws.Range("H17").Interior.Color = 13434828

# ---------------------------------------------------------------------------
# 5) Append two brand-new firms that just initiated coverage on 6/26/2018.
#    Like all other rows, only columns A-D carry data (E:H are left blank,
#    matching the existing sheet convention for newly added rows).
# ---------------------------------------------------------------------------
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "6/26/2018,Initiates,Buy,`$280.00"
$ws.Range("C28").Value = "6/26/2018,Initiates,Buy,`$280.00"
$ws.Range("D28").Value = "6/26/2018,Initiates,Buy,`$280.00"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
